$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 145 (shifts old rows 145..162 down to 146..163)
$ws.Rows.Item(145).Insert()

# Populate the newly inserted row 145 with the new data record
$ws.Range("A145").Value = 4
$ws.Range("B145").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C145").Value = "Los Lagos"
$ws.Range("D145").Value = 45142
$ws.Range("E145").Value = 10
$ws.Range("F145").Value = 100112031
$ws.Range("G145").Value = "Poroto verde"
$ws.Range("H145").Value = "Magnum"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 45
$ws.Range("K145").Value = 30000
$ws.Range("L145").Value = 30000
$ws.Range("M145").Value = 30000
$ws.Range("N145").Value = "`$/malla 25 kilos"
$ws.Range("O145").Value = "Perú"
$ws.Range("P145").Value = 1200
$ws.Range("Q145").Value = 25
$ws.Range("R145").Value = "Hortaliza"
